$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (45572 = 2024-10-07).
# Update every data row (2 through 29) to the new serial value 45573 (2024-10-08),
# keeping the existing numeric storage/format (style index 1, numFmt yyyy-mm-dd).
for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 3).Value = 45573
}
